$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must remain TEXT (many look numeric, e.g. "1.75", or
# have two dots like "35.374.11" and are not valid numbers at all). Excel's COM .Value
# setter auto-converts plain-numeric-looking strings to real numbers, so we temporarily
# force the whole D2:D51 range to Text format, write every value as a string, then clear
# the temporary formatting again so no visible style/number-format change remains.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "35.374.11"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "1.890.47"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("E4").Value = "  -0.79%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "246.79"
$ws.Range("E5").Value = "  -2.27%  "

$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "0.692"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("D8").Value = "43.35"
$ws.Range("E8").Value = "  +5.30%  "

$ws.Range("E9").Value = "  -1.94%  "

$ws.Range("D10").Value = "53.78"
$ws.Range("E10").Value = "  +1.66%  "

$ws.Range("D11").Value = "0.0742"
$ws.Range("E11").Value = "  -1.75%  "

$ws.Range("D13").Value = "13.37"
$ws.Range("E13").Value = "  +2.72%  "

$ws.Range("D14").Value = "2.165.65"
$ws.Range("E14").Value = "  -0.60%  "

$ws.Range("D15").Value = "0.761"
$ws.Range("E15").Value = "  +3.03%  "

$ws.Range("D16").Value = "4.91"
$ws.Range("E16").Value = "  -1.27%  "

$ws.Range("D17").Value = "1.891.45"
$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").Value = "35.447.02"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("D19").Value = "73.21"
$ws.Range("E19").Value = "  -0.94%  "

$ws.Range("E20").Value = "  -1.29%  "

$ws.Range("D21").Value = "245.15"
$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("E23").Value = "  -1.78%  "

$ws.Range("D24").Value = "2.67"
$ws.Range("E24").Value = "  +9.82%  "

$ws.Range("E25").Value = "  -0.76%  "

$ws.Range("E26").Value = "  -6.53%  "

$ws.Range("D27").Value = "166.09"
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("E28").Value = "  -1.22%  "

$ws.Range("E29").Value = "  -0.60%  "

$ws.Range("E30").Value = "  -1.66%  "

$ws.Range("D31").Value = "4.128.44"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").Value = "1.75"
$ws.Range("E32").Value = "  +11.44%  "

$ws.Range("E33").Value = "  -1.30%  "

$ws.Range("D34").Value = "0.0586"
$ws.Range("E34").Value = "  -4.19%  "

$ws.Range("D35").Value = "4.18"
$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("E36").Value = "  -13.04%  "

$ws.Range("E37").Value = "  -0.91%  "

$ws.Range("D38").Value = "0.846"
$ws.Range("E38").Value = "  -1.41%  "

$ws.Range("E39").Value = "  -1.76%  "

$ws.Range("E40").Value = "  +6.70%  "

$ws.Range("E41").Value = "  +2.59%  "

$ws.Range("D42").Value = "17.32"
$ws.Range("E42").Value = "  +0.81%  "

$ws.Range("D43").Value = "97.04"
$ws.Range("E43").Value = "  -2.72%  "

$ws.Range("E44").Value = "  -2.77%  "

$ws.Range("D45").Value = "1.297.80"
$ws.Range("E45").Value = "  -2.57%  "

$ws.Range("D46").Value = "2.33"
$ws.Range("E46").Value = "  -4.34%  "

$ws.Range("D47").Value = "0.0796"
$ws.Range("E47").Value = "  +7.23%  "

$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").Value = "12.34"
$ws.Range("E49").Value = "  +4.79%  "

$ws.Range("E50").Value = "  -0.61%  "

$ws.Range("D51").Value = "6.26"
$ws.Range("E51").Value = "  -5.42%  "

# Remove the temporary text-number-format now that every price has been written as a string;
# this restores the original (default/General) formatting on the price column.
$priceRange.ClearFormats()
